$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list: refresh Price (D) and Volume(1h) (E) columns.
# Values are kept as text (matching the original inlineStr cells), so force
# a text number format before assigning to avoid Excel auto-converting them
# to numeric/percentage values.
$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","D12","E12","E13","D14","E14","D15","E15","D16","E16","E17","D18","E18","E19","E20","D21","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","D27","D39","E39","D40","E40","D41","E41","E43","D44","E44","D45","E45","D46","E46","E47","E48","D49","E49","E50","E51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "290.54"
$ws.Range("E2").Value = "-1.02%"
$ws.Range("D3").Value = "30.81"
$ws.Range("E3").Value = "-1.52%"
$ws.Range("D4").Value = "4.885"
$ws.Range("E4").Value = "-1.73%"
$ws.Range("D5").Value = "0.07237"
$ws.Range("E5").Value = "-1.51%"
$ws.Range("D6").Value = "2.362"
$ws.Range("E6").Value = "29.99%"
$ws.Range("D7").Value = "7.661"
$ws.Range("E7").Value = "-0.01%"
$ws.Range("D8").Value = "3.711"
$ws.Range("E8").Value = "-1.12%"
$ws.Range("D9").Value = "0.8973"
$ws.Range("E9").Value = "-1.29%"
$ws.Range("D10").Value = "0.1665"
$ws.Range("E10").Value = "1.04%"
$ws.Range("D11").Value = "0.07946"
$ws.Range("D12").Value = "0.08150"
$ws.Range("E12").Value = "-0.28%"
$ws.Range("E13").Value = "3.70%"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").Value = "0.48%"
$ws.Range("D15").Value = "0.001497"
$ws.Range("E15").Value = "-0.01%"
$ws.Range("D16").Value = "0.005833"
$ws.Range("E16").Value = "3.45%"
$ws.Range("E17").Value = "0.44%"
$ws.Range("D18").Value = "2.074"
$ws.Range("E18").Value = "-2.21%"
$ws.Range("E19").Value = "0.99%"
$ws.Range("E20").Value = "-0.76%"
$ws.Range("D21").Value = "3.981"
$ws.Range("E21").Value = "-8.18%"
$ws.Range("D22").Value = "0.2305"
$ws.Range("E22").Value = "16.64%"
$ws.Range("D23").Value = "0.04519"
$ws.Range("E23").Value = "0.98%"
$ws.Range("E24").Value = "-1.20%"
$ws.Range("D25").Value = "0.004394"
$ws.Range("E25").Value = "8.43%"
$ws.Range("D26").Value = "0.0001303"
$ws.Range("E26").Value = "4.04%"
$ws.Range("D27").Value = "0.0003399"
$ws.Range("D39").Value = "0.01576"
$ws.Range("E39").Value = "-3.84%"
$ws.Range("D40").Value = "0.04374"
$ws.Range("E40").Value = "-0.54%"
$ws.Range("D41").Value = "0.007244"
$ws.Range("E41").Value = "-2.67%"
$ws.Range("E43").Value = "-0.99%"
$ws.Range("D44").Value = "0.002027"
$ws.Range("E44").Value = "-1.51%"
$ws.Range("D45").Value = "0.009506"
$ws.Range("E45").Value = "-14.51%"
$ws.Range("D46").Value = "0.00005726"
$ws.Range("E46").Value = "-4.21%"
$ws.Range("E47").Value = "0.08%"
$ws.Range("E48").Value = "12.69%"
$ws.Range("D49").Value = "0.002902"
$ws.Range("E49").Value = "-3.28%"
$ws.Range("E50").Value = "0.08%"
$ws.Range("E51").Value = "0.08%"
